$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Insert()
$ws.Range("A7").Value = "TSM"

$ws.Range("A14").Insert()
$ws.Range("A14").Value = "MU"

$ws.Range("A34").Insert()
$ws.Range("A34").Value = "RTX"

$ws.Range("A40").Insert()
$ws.Range("A40").Value = "SHEL"

$ws.Range("A73:A74").Insert()
$ws.Range("A73").Value = "VRTX"
$ws.Range("A74").Value = "PGR"

$ws.Range("A85").Insert()
$ws.Range("A85").Value = "ICE"

$ws.Range("A114").Insert()
$ws.Range("A114").Value = "GS-PA"

$ws.Range("A128").Insert()
$ws.Range("A128").Value = "FTNT"

$ws.Range("A146:A147").Insert()
$ws.Range("A146").Value = "F"
$ws.Range("A147").Value = "CMG"

$ws.Range("A155").Insert()
$ws.Range("A155").Value = "CBRE"

$ws.Range("A191").Insert()
$ws.Range("A191").Value = "PEG"

$ws.Range("A209").Insert()
$ws.Range("A209").Value = "XYL"

$ws.Range("A228:A229").Insert()
$ws.Range("A228").Value = "COHR"
$ws.Range("A229").Value = "EXR"

$ws.Range("A240:A241").Insert()
$ws.Range("A240").Value = "LPLA"
$ws.Range("A241").Value = "UMC"

$ws.Range("A243").Insert()
$ws.Range("A243").Value = "NTRS"

$ws.Range("A275").Insert()
$ws.Range("A275").Value = "VRSN"

$ws.Range("A321").Insert()
$ws.Range("A321").Value = "PHYS"

$ws.Range("A323").Insert()
$ws.Range("A323").Value = "SMCI"

$ws.Range("A335").Insert()
$ws.Range("A335").Value = "ENTG"

$ws.Range("A340").Insert()
$ws.Range("A340").Value = "ZBH"

$ws.Range("A351").Insert()
$ws.Range("A351").Value = "WES"

$ws.Range("A370").Insert()
$ws.Range("A370").Value = "FFIV"

$ws.Range("A394").Insert()
$ws.Range("A394").Value = "TSEM"

$ws.Range("A421").Insert()
$ws.Range("A421").Value = "UNM"

$ws.Range("A477").Insert()
$ws.Range("A477").Value = "XP"

$ws.Range("A481").Insert()
$ws.Range("A481").Value = "ARMK"

$ws.Range("A492").Insert()
$ws.Range("A492").Value = "VNO-PM"

$ws.Range("A496").Insert()
$ws.Range("A496").Value = "TTMI"

$ws.Range("A510").Delete()
